$d = $word.ActiveDocument

$xml25 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>shift_table</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> will contain an entry for each shift. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Containing  columns</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>, linking to the user and storing the date. An entry will be created for each and every shift</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>shiftID</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> is a big </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, as there will be 5 entries per week for each member of staff. So 6 staff * 5 days is 30 shifts per week, or 1,560 per year. Having </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>BigInt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> allows for expansion of the system, and ensures longevity. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Stored Procedures</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>User_add</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> --</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>User_get</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> --</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>User_getAll</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> --</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>User_delete</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> --</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>User_edit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> --</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Level_add</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> --</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Level_delete</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> --</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Level_edit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> --</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Level_get</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> --</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Level_getAll</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> --</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Shift_add</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> --</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Shift_delete</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> --</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Shift_edit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> --</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Shift_get</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t xml:space="preserve">Gets a shift based on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>shiftID</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:lastRenderedPageBreak/><w:t>Shift_getStartDate</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t>Get all the shifts for a week, based on the start date passed in</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Shift_getUser</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t xml:space="preserve">Get all the shifts for a user, based on </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$p = $d.Paragraphs(25)
$p.Range.InsertXML($xml25)

$xml23 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>shift</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_date</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:tab/><w:t>date</w:t></w:r></w:p>
'@
$p = $d.Paragraphs(23)
$p.Range.InsertXML($xml23)

$xml22 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>userID</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:t>FK</w:t></w:r></w:p>
'@
$p = $d.Paragraphs(22)
$p.Range.InsertXML($xml22)

$xml21 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>shiftID</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Big</w:t></w:r><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:t>Identity</w:t></w:r><w:r><w:tab/><w:t>PK</w:t></w:r></w:p>
'@
$p = $d.Paragraphs(21)
$p.Range.InsertXML($xml21)

$xml20 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Shift_table</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$p = $d.Paragraphs(20)
$p.Range.InsertXML($xml20)

$xml18 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>level_table</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> will contain each level of nurse</w:t></w:r></w:p>
'@
$p = $d.Paragraphs(18)
$p.Range.InsertXML($xml18)

$xml16 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>level</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:tab/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nvachar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(100)</w:t></w:r></w:p>
'@
$p = $d.Paragraphs(16)
$p.Range.InsertXML($xml16)

$xml15 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>levelID</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:tab/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:t>FK</w:t></w:r></w:p>
'@
$p = $d.Paragraphs(15)
$p.Range.InsertXML($xml15)

$xml14 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Level_table</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$p = $d.Paragraphs(14)
$p.Range.InsertXML($xml14)

$xml12 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_table</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> will contain each user that has access to the system. </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>levelID</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> is a foreign key, referencing the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>level_</w:t></w:r><w:r><w:t>table</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p>
'@
$p = $d.Paragraphs(12)
$p.Range.InsertXML($xml12)

$xml10 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>staffID</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$p = $d.Paragraphs(10)
$p.Range.InsertXML($xml10)

$xml9 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>levelID</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:tab/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>FK</w:t></w:r></w:p>
'@
$p = $d.Paragraphs(9)
$p.Range.InsertXML($xml9)

$xml8 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>password</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:tab/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nvarchar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(100)</w:t></w:r></w:p>
'@
$p = $d.Paragraphs(8)
$p.Range.InsertXML($xml8)

$xml7 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>forename</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:tab/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nvarchar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(200)</w:t></w:r></w:p>
'@
$p = $d.Paragraphs(7)
$p.Range.InsertXML($xml7)

$xml6 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>surname</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:tab/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nvarchar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(200)</w:t></w:r><w:r><w:tab/></w:r></w:p>
'@
$p = $d.Paragraphs(6)
$p.Range.InsertXML($xml6)

$xml5 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>userID</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Identity</w:t></w:r><w:r><w:tab/><w:t>PK</w:t></w:r></w:p>
'@
$p = $d.Paragraphs(5)
$p.Range.InsertXML($xml5)

$xml4 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>User_table</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$p = $d.Paragraphs(4)
$p.Range.InsertXML($xml4)
